$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink column A width from 23 to 12 (stored OOXML width = ColumnWidth + 5/6)
$ws.Columns.Item(1).ColumnWidth = 11.166666666666666

# Trim car names down to just the brand/make
$ws.Cells.Item(2, 1).Value = "BYD"
$ws.Cells.Item(3, 1).Value = "XPENG"
$ws.Cells.Item(4, 1).Value = "VINFAST"
$ws.Cells.Item(5, 1).Value = "Honda"
$ws.Cells.Item(6, 1).Value = "BYD"
$ws.Cells.Item(7, 1).Value = "Volkswagen"
$ws.Cells.Item(8, 1).Value = "BMW"
$ws.Cells.Item(9, 1).Value = "smart"
$ws.Cells.Item(10, 1).Value = "BYD"
$ws.Cells.Item(11, 1).Value = "Hyundai"
$ws.Cells.Item(12, 1).Value = "Kia"
$ws.Cells.Item(13, 1).Value = "NIO"
$ws.Cells.Item(14, 1).Value = "NIO"
$ws.Cells.Item(15, 1).Value = "Lexus"
